# Add 2022-Q4 data.
#
# Before: 总计 , 2022-Q2 (fund holdings for 2022-Q2)
# After : 总计 , 2022-Q4 (fund holdings for 2022-Q4, NEW) , 2022-Q2 (unchanged, moved)
#
# The existing "2022-Q2" worksheet (sheetId 2) is duplicated so the old
# data survives unmodified under the same name (now sheetId 3, 3rd tab).
# The original worksheet object is renamed to "2022-Q4" and its contents
# are replaced with the Q4 fund holdings (still sheetId 2, 2nd tab).
# Finally the "总计" (totals) sheet gets a new row for 2022-Q4.

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item(1)
$q2 = $wb.Worksheets.Item(2)

# 1) Duplicate the current "2022-Q2" sheet right after itself so its data
#    is preserved verbatim under the old name.
$q2.Copy($null, $q2)
$q2Copy = $wb.Worksheets.Item(3)

# 2) Free up the "2022-Q2" name on the original sheet and rename it to
#    "2022-Q4" - this keeps sheetId 2 on the Q4 data and gives the
#    untouched duplicate sheetId 3.
$q2.Name = "2022-Q4"
$q2Copy.Name = "2022-Q2"
$q4 = $q2

# 3) Wipe the old Q2 fund rows out of the (renamed) Q4 sheet and rebuild
#    it with the Q4 fund holdings.
$q4.Range("A1:H3").Clear()

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$q4.Range("A2").Value = 0

$q4.Range("B2").NumberFormat = "@"
$q4.Range("B2").Value = "015245"
$q4.Range("C2").Value = "南华丰汇混合"

$q4.Range("D2").NumberFormat = "@"
$q4.Range("D2").Value = "0.11"
$q4.Range("E2").NumberFormat = "@"
$q4.Range("E2").Value = "84.24"
$q4.Range("F2").NumberFormat = "@"
$q4.Range("F2").Value = "1.08"
$q4.Range("G2").NumberFormat = "@"
$q4.Range("G2").Value = "0.0012"

$q4.Range("H2").Value = 5

# Match styling: the new header row / index column reuse the "总计"
# sheet's header style, and page margins follow the same sheet.
$total.Range("B1:D1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$total.Range("A2").Copy()
$q4.Range("A2").PasteSpecial(-4122)

# The text-valued cells shouldn't carry the "@" number-format style once
# their (already-text) value is committed.
$q4.Range("B2").ClearFormats()
$q4.Range("D2:G2").ClearFormats()

$q4.PageSetup.LeftMargin = 54
$q4.PageSetup.RightMargin = 54
$q4.PageSetup.TopMargin = 72
$q4.PageSetup.BottomMargin = 72
$q4.PageSetup.HeaderMargin = 36
$q4.PageSetup.FooterMargin = 36

# 4) Insert the new "2022-Q4" row into the "总计" summary sheet, shifting
#    the existing "2022-Q2" row down.
$total.Rows(2).Insert()
$total.Range("B2:D2").ClearFormats()

$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0

$total.Range("A3").Value = 1
